$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12)
# for every data row (2 through 288). Update them all in one go.
$ws.Range("C2:C288").Value = 45181
